$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B30 with new value
$ws.Range("B30").Value = 0.39639999999999997

# Recalculate so AVERAGE formula in B32 updates
$excel.Calculate()

# Move the active selection from B28 to C1
$ws.Range("C1").Select()
